# Update "Correspond Handback/Handoff" and "Latest HO Xliff Generate Date" timestamps
# in the handback status report, as part of "Generate Report for Handback".

$wb = $excel.ActiveWorkbook

# Overview sheet: Latest HO Xliff Generate Date for the first row (d1ac55cc...)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-27 21:03:18"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for first row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-27 21:03:14"
$wsZhCn.Range("K2").Value = "2016-08-27 21:03:31"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for first row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-27 21:03:18"
$wsDeDe.Range("K2").Value = "2016-08-27 21:03:38"
